$wb = $excel.ActiveWorkbook

# Add the Projects worksheet after the last existing sheet (Clients)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$projects = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$projects.Name = "Projects"

# Add the Tasks worksheet after Projects
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$tasks = $wb.Worksheets.Add([Type]::Missing, $lastSheet2)
$tasks.Name = "Tasks"

# Populate Projects headers
$projects.Range("A1").Value = "TITLE"
$projects.Range("B1").Value = "CLIENT"
$projects.Range("C1").Value = "HOUR"
$projects.Range("D1").Value = "PRIORITY"
$projects.Range("E1").Value = "START_DATE"
$projects.Range("F1").Value = "END_DATE"
$projects.Range("G1").Value = "SUMMARY"
$projects.Range("H1").Value = "TEAM"
$projects.Range("I1").Value = "DESCRIPTION"

# Populate Tasks headers
$tasks.Range("A1").Value = "TITLE"
$tasks.Range("B1").Value = "START_DATE"
$tasks.Range("C1").Value = "END_DATE"
$tasks.Range("D1").Value = "HOUR"
$tasks.Range("E1").Value = "PROJECT"
$tasks.Range("F1").Value = "SUMMARY"
$tasks.Range("G1").Value = "DESCRIPTION"
$tasks.Range("H1").Value = "STATUS"
$tasks.Range("I1").Value = "ATTACH"

# Update the selection on the Login sheet to span full columns F:XFD
$login = $wb.Worksheets.Item("Login")
$login.Range("F1:XFD1048576").Select() | Out-Null

# Remove the row-level custom format flag on row 1 of Login while keeping
# the individual cell styles intact.
$row1 = $login.Rows.Item(1)
$row1.ClearFormats()
$headerRange = $login.Range("A1:E1")
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$headerRange.Interior.Color = 65535

# Restore the per-sheet selections that were active before saving.
$projects.Activate()
$projects.Range("C1:G1").Select() | Out-Null

$tasks.Activate()
$tasks.Range("J2").Select() | Out-Null
